# Auto-generated edit script applying crypto price/volume/coin updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '90.227.09'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +2.40%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.212.05'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.65%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '214.83'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.24%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '619.50'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.390'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.701'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '3.205.59'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.67%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.579'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.70%  '
$ws.Range('E12').Value = '  -3.98%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000259'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('B14').Value = 'Toncoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.44'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.25%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.812.47'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.59%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '90.121.79'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.41%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '32.77'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.77%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.228.56'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.00%  '
$ws.Range('B19').Value = 'PEPE'
$ws.Range('C19').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0000237'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +69.37%  '
$ws.Range('B20').Value = 'SuiNetwork'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.34'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +8.84%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '440.01'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.41%  '
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '13.39'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.80%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.60'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -3.20%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.05'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.05%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.09'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.74%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.52'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -4.60%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '3.381.13'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.92%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '75.17'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.95%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.161'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -12.49%  '
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.17'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +33.10%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '8.43'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.71%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '535.23'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.85%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.86'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.86%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.89'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.18%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.26'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -4.77%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '22.48'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.46%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '22.34'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.31%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('E41').Value = '  -8.29%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('E43').Value = '  -2.21%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.374'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.75%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '150.95'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '171.71'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.96%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '43.32'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.48%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.124'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -6.29%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.738'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.21%  '
$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.23'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -7.13%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.617'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.20%  '
